# Cryptos list refresh (GitHub Actions daily update).
# All Price (D) / Volume(1h) (E) cells are stored as plain text in this
# sheet, so any value that Excel would otherwise auto-parse as a number
# gets NumberFormat "@" applied first to keep it text (matches the
# original t="inlineStr" cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.707.64"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "3.336.51"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.15"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.27"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "3.334.58"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.179"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.578"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.43"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "673.90"
$ws.Range("E14").Value = "  +5.46%  "
$ws.Range("D15").Value = "3.877.24"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.41"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "67.762.38"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "3.330.66"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.40"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.889"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.45"
$ws.Range("E23").Value = "  +9.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.12"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.06"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.78"
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.44"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("E31").Value = "  +10.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "572.45"
$ws.Range("E32").Value = "  -3.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.98"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "3.693.85"
$ws.Range("E36").Value = "  -6.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.65"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("E38").Value = "  -5.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.49"
$ws.Range("E39").Value = "  +5.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.333"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0667"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0406"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.61"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.08"
$ws.Range("E51").Value = "  -1.21%  "
